$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("Ja Morant", "PG", "Memphis Grizzlies"),
    @("Scottie Barnes", "SG,SF,PF", "Toronto Raptors"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("Guerschon Yabusele", "PF,C", "Philadelphia 76ers"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
